$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.849.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.652"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.10"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.061.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "602.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.995.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.505.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.44%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.55"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +27.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.17"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.728.80"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.51%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.80"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "498.31"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.49%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0456"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.63"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000243"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.95"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.91%  "
